$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / non-numeric-looking values: direct assignment
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("E20").Value = "19BitKanKAN"

# Numeric-looking values that must remain text: force text entry via leading apostrophe
$ws.Range("D2").Value = "'242.91"
$ws.Range("D3").Value = "'23.00"
$ws.Range("D4").Value = "'5.401"
$ws.Range("D5").Value = "'0.05900"
$ws.Range("D6").Value = "'3.453"
$ws.Range("D7").Value = "'6.587"
$ws.Range("D8").Value = "'0.8107"
$ws.Range("D9").Value = "'0.9182"
$ws.Range("D10").Value = "'0.01132"
$ws.Range("D11").Value = "'0.1410"
$ws.Range("D12").Value = "'0.07420"
$ws.Range("D13").Value = "'0.03270"
$ws.Range("D14").Value = "'0.03064"
$ws.Range("D15").Value = "'0.09343"
$ws.Range("D16").Value = "'3.868"
$ws.Range("D17").Value = "'0.001557"
$ws.Range("D18").Value = "'0.04677"
$ws.Range("D19").Value = "'0.005892"
$ws.Range("D20").Value = "'0.001269"
$ws.Range("D21").Value = "'0.004900"
$ws.Range("D22").Value = "'0.00009003"
$ws.Range("D24").Value = "'2.143"
$ws.Range("D40").Value = "'0.03953"
$ws.Range("D44").Value = "'0.009296"
$ws.Range("D45").Value = "'0.00005205"
$ws.Range("D47").Value = "'0.7503"
$ws.Range("D48").Value = "'0.002285"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002001"
